# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" (column E) values for rows 16-28 are reversed in order
# (the oldest periods become the newest ones and vice versa), and the
# "Valor Mora" (column F) values on the first and last data rows are swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Periodo Mora" values for rows 16..28 (reversed order of the old list)
$periodos = @("2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002")

$row = 16
foreach ($p in $periodos) {
    $ws.Range("E$row").Value = $p
    $row = $row + 1
}

# Swap the "Valor Mora" values between the first (row 16) and last (row 28) data rows
$ws.Range("F16").Value = 25749
$ws.Range("F28").Value = 35112
